$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1) with the new/renamed column labels.
$ws.Range("A1").Value = "Z_real"
$ws.Range("B1").Value = "Z_imag"
$ws.Range("C1").Value = "angular frequency"
$ws.Range("D1").Value = "eff_cap"
$ws.Range("E1").Value = "applied voltage"
$ws.Range("F1").Value = "J_ph"
$ws.Range("G1").Value = "J"

# Move/restore the active selection to B2, matching the saved view state.
$ws.Range("B2").Select()
